$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$s.Shapes.Item("TextBox 7").Delete()
